# Fabric + DepencyInhection + logger
# Update unemployment report data: refresh Ultimo/Anterior values and
# Referencia periods for several countries, and fix the row ordering for
# a few country pairs whose relative ranking changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - Peru
$ws.Range("B7").Value = 93.2
$ws.Range("C7").Value = 92.8
$ws.Range("D7").Value = "2022-06"

# Row 10 - Holanda
$ws.Range("B10").Value = 81.2
$ws.Range("C10").Value = 81
$ws.Range("D10").Value = "2022-03"

# Row 11 - Suíça
$ws.Range("B11").Value = 79.3
$ws.Range("C11").Value = 80.40000000000001
$ws.Range("D11").Value = "2022-03"

# Row 13 - Alemanha
$ws.Range("B13").Value = 77
$ws.Range("C13").Value = 76.90000000000001
$ws.Range("D13").Value = "2022-03"

# Row 14 - Malta
$ws.Range("B14").Value = 76.09999999999999
$ws.Range("C14").Value = 76.7
$ws.Range("D14").Value = "2022-03"

# Row 16 - was "República Tcheca", now "Bolívia"
$ws.Range("A16").Value = "Bolívia"
$ws.Range("B16").Value = 75.2
$ws.Range("C16").Value = 73.59999999999999
# D16 stays "2021-12"

# Row 17 - was "Bolívia", now "República Tcheca"
$ws.Range("A17").Value = "República Tcheca"
$ws.Range("B17").Value = 75
$ws.Range("C17").Value = 75.3
$ws.Range("D17").Value = "2022-03"

# Row 19 - Áustria
$ws.Range("B19").Value = 73.3
$ws.Range("C19").Value = 73.5
$ws.Range("D19").Value = "2022-03"

# Row 20 - was "Lituânia", now "Irlanda"
$ws.Range("A20").Value = "Irlanda"
$ws.Range("B20").Value = 72.8
$ws.Range("C20").Value = 73
$ws.Range("D20").Value = "2022-03"

# Row 21 - was "Irlanda", now "Lituânia"
$ws.Range("A21").Value = "Lituânia"
$ws.Range("B21").Value = 72.59999999999999
$ws.Range("C21").Value = 73.2
# D21 stays "2022-03"

# Row 22 - Eslovenia
$ws.Range("B22").Value = 72.5
$ws.Range("C22").Value = 72.40000000000001
$ws.Range("D22").Value = "2022-03"

# Row 23 - Chipre
$ws.Range("B23").Value = 72
# C23 stays 72.3
$ws.Range("D23").Value = "2022-03"

# Row 24 - Polônia (values unchanged, only Referência date changes)
$ws.Range("D24").Value = "2022-03"

# Row 25 - was "Eslováquia", now "Luxemburgo"
$ws.Range("A25").Value = "Luxemburgo"
$ws.Range("B25").Value = 70.59999999999999
$ws.Range("C25").Value = 69.40000000000001
$ws.Range("D25").Value = "2022-03"

# Row 26 - was "Noruega", now "Eslováquia"
$ws.Range("A26").Value = "Eslováquia"
$ws.Range("B26").Value = 70.59999999999999
$ws.Range("C26").Value = 70.8
# D26 stays "2022-03"

# Row 27 - was "Luxemburgo", now "Noruega"
$ws.Range("A27").Value = "Noruega"
$ws.Range("B27").Value = 69.90000000000001
$ws.Range("C27").Value = 69.8
$ws.Range("D27").Value = "2022-03"

# Row 34 - Bélgica
$ws.Range("B34").Value = 66.40000000000001
$ws.Range("C34").Value = 66.2
$ws.Range("D34").Value = "2022-03"

# Row 35 - Croácia
$ws.Range("B35").Value = 64.2
$ws.Range("C35").Value = 64.09999999999999
$ws.Range("D35").Value = "2022-03"
